$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (revision bump: "(4)" -> "(5)")
$ws.Name = "repayment_20250912_20250912 (5)"

# --- Numeric (Talk_time / Repayment_collections) column updates ---
$ws.Range("H2").Value  = 348
$ws.Range("H3").Value  = 1.8280000000000001
$ws.Range("H4").Value  = 2.1429999999999998
$ws.Range("H5").Value  = 1.3759999999999999
$ws.Range("H6").Value  = 1.385

$ws.Range("D7").Value  = 4
$ws.Range("H7").Value  = 1.544

$ws.Range("H8").Value  = 1.165
$ws.Range("H9").Value  = 785
$ws.Range("H10").Value = 1.464

$ws.Range("D11").Value = 1
$ws.Range("H11").Value = 1.167

$ws.Range("H12").Value = 1.3009999999999999

$ws.Range("D15").Value = 5
$ws.Range("H15").Value = 4.9379999999999997

$ws.Range("H16").Value = 1.016

$ws.Range("D17").Value = 1
$ws.Range("H17").Value = 1.494

$ws.Range("H18").Value = 924

# --- Text columns (Repayment_amount / Pending Amount Recovery) ---
# these hold numeric-looking text, so force a text number-format before
# assigning, otherwise Excel would coerce "1,559,407.00" -> the number
# 1559407 instead of keeping it as a literal string.
function Set-TextValue($rng, $text) {
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

Set-TextValue $ws.Range("E7")  "3,904,049.00"
Set-TextValue $ws.Range("G7")  "2.05"

Set-TextValue $ws.Range("E11") "1,559,407.00"
Set-TextValue $ws.Range("G11") "1.05"

Set-TextValue $ws.Range("E15") "2,996,580.00"
Set-TextValue $ws.Range("G15") "1.97"

Set-TextValue $ws.Range("E17") "89,732.00"
Set-TextValue $ws.Range("G17") "0.06"

# --- Selection change ---
$ws.Range("A3").Select()
